$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Standard Deviation values updated
$ws.Range("B3").Value = 0.02588004227551893
$ws.Range("C3").Value = 0.01525519540867466
$ws.Range("D3").Value = 0.04415322086835657
$ws.Range("E3").Value = 0.03219527645924017
$ws.Range("F3").Value = 0.07235739080727961
$ws.Range("G3").Value = 0.1071726775388513

# Row 4: Maximum - only B4 changes
$ws.Range("B4").Value = 0.07724121144269158
